$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO" (columns A:R, data rows 2..6, totals row 7)
# A new client "PAREDES POVEDA TATIANA VERONICA" (all-zero sales) is
# inserted in alphabetical order right before "VACA PANCHI DORYS
# CAROLINA", pushing it (and the totals row) down by one row. The
# totals row's "X de 5" labels become "X de 6" to reflect the new
# total row count.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row 6; this shifts old row 6 (VACA PANCHI DORYS CAROLINA)
# to row 7 and the totals row from 7 to 8, carrying their values/styles.
$ws1.Rows.Item(6).Insert()

$ws1.Cells.Item(6, 1).Value = "VACA PANCHI CAROLINA"
$ws1.Cells.Item(6, 2).Value = "PAREDES POVEDA TATIANA VERONICA"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(6, $col).Value = 0
}

# Fix up the totals row (now row 8): "X de 5" -> "X de 6"
$ws1.Cells.Item(8, 3).Value  = "0 de 6"
$ws1.Cells.Item(8, 4).Value  = "1 de 6"
$ws1.Cells.Item(8, 5).Value  = "0 de 6"
$ws1.Cells.Item(8, 6).Value  = "0 de 6"
$ws1.Cells.Item(8, 7).Value  = "0 de 6"
$ws1.Cells.Item(8, 8).Value  = "0 de 6"
$ws1.Cells.Item(8, 9).Value  = "0 de 6"
$ws1.Cells.Item(8, 10).Value = "0 de 6"
$ws1.Cells.Item(8, 11).Value = "0 de 6"
$ws1.Cells.Item(8, 12).Value = "0 de 6"
$ws1.Cells.Item(8, 13).Value = "2 de 6"
$ws1.Cells.Item(8, 14).Value = "0 de 6"
$ws1.Cells.Item(8, 15).Value = "0 de 6"
$ws1.Cells.Item(8, 16).Value = "0 de 6"
$ws1.Cells.Item(8, 17).Value = "0 de 6"
$ws1.Cells.Item(8, 18).Value = "0 de 6"

# ------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL" (columns A:G, data rows 2..6, totals row 7)
# Two new all-zero clients are inserted in alphabetical order:
#   "KITCHENSCO S.A."                 before LINCANGO LUGMANIA SANDY LIZETH
#   "PAREDES POVEDA TATIANA VERONICA" before VACA PANCHI DORYS CAROLINA
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Insert new row 5 (KITCHENSCO S.A.); old row 5 (LINCANGO...) becomes row 6,
# old row 6 (VACA PANCHI DORYS CAROLINA) becomes row 7, totals row becomes 8.
$ws2.Rows.Item(5).Insert()
$ws2.Cells.Item(5, 1).Value = "VACA PANCHI CAROLINA"
$ws2.Cells.Item(5, 2).Value = "KITCHENSCO S.A."
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(5, $col).Value = 0
}

# Insert new row 7 (PAREDES POVEDA TATIANA VERONICA); old row 7
# (VACA PANCHI DORYS CAROLINA, after the previous shift) becomes row 8,
# and the totals row becomes row 9.
$ws2.Rows.Item(7).Insert()
$ws2.Cells.Item(7, 1).Value = "VACA PANCHI CAROLINA"
$ws2.Cells.Item(7, 2).Value = "PAREDES POVEDA TATIANA VERONICA"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(7, $col).Value = 0
}
